$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '57.128.09'
$ws.Range('E2').Value = '  -3.05%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.549.17'
$ws.Range('E3').Value = '  -4.25%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '516.90'
$ws.Range('E5').Value = '  -1.19%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '141.28'
$ws.Range('E6').Value = '  -2.11%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.16%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.558'
$ws.Range('E8').Value = '  -2.03%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.552.70'
$ws.Range('E9').Value = '  -4.50%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.58'
$ws.Range('E10').Value = '  -5.57%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0990'
$ws.Range('E11').Value = '  -3.56%  '
$ws.Range('E12').Value = '  -3.67%  '
$ws.Range('E13').Value = '  -0.29%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.985.49'
$ws.Range('E14').Value = '  -4.72%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '57.137.66'
$ws.Range('E15').Value = '  -3.05%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '20.01'
$ws.Range('E16').Value = '  -4.81%  '
$ws.Range('E17').Value = '  -3.10%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.509.75'
$ws.Range('E18').Value = '  -6.44%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '330.70'
$ws.Range('E19').Value = '  -2.51%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '4.25'
$ws.Range('E20').Value = '  -3.28%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '10.08'
$ws.Range('E21').Value = '  -2.71%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.14'
$ws.Range('E22').Value = '  -3.78%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '64.76'
$ws.Range('E24').Value = '  +0.60%  '
$ws.Range('E25').Value = '  +0.60%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.997'
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.399'
$ws.Range('E27').Value = '  -4.91%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.651.26'
$ws.Range('E28').Value = '  -4.86%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.89'
$ws.Range('E29').Value = '  -3.62%  '
$ws.Range('E30').Value = '  -0.06%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.0₃0737'
$ws.Range('E31').Value = '  -7.99%  '
$ws.Range('E32').Value = '  -6.32%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.57'
$ws.Range('E33').Value = '  -1.80%  '
$ws.Range('E34').Value = '  -1.43%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '18.44'
$ws.Range('E35').Value = '  -2.40%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.97'
$ws.Range('E36').Value = '  -4.33%  '
$ws.Range('E37').Value = '  -4.76%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.829'
$ws.Range('E38').Value = '  -7.79%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '35.58'
$ws.Range('E39').Value = '  -3.46%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.823'
$ws.Range('E40').Value = '  -5.49%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.43'
$ws.Range('E41').Value = '  -2.41%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.998'
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '3.46'
$ws.Range('E43').Value = '  -3.45%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '10.63'
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('E45').Value = '  -2.26%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '264.96'
$ws.Range('E46').Value = '  -3.82%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.577'
$ws.Range('E47').Value = '  -6.29%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '18.52'
$ws.Range('E48').Value = '  -6.53%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0514'
$ws.Range('E49').Value = '  -3.55%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.954.55'
$ws.Range('E50').Value = '  -4.71%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '4.51'
$ws.Range('E51').Value = '  -4.61%  '
